# Auto-generated script applying scheduled market-data refresh to Sheets/Phoenix_Profits.xlsx
# Updates currentAveragePrice* / Leve price & profit columns (H:N) for specific rows per sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3451.889
$ws.Range("I53").Value = 6821.875
$ws.Range("J53").Value = 755.9
$ws.Range("K53").Value = 6821.875
$ws.Range("L53").Value = 755.9
$ws.Range("M53").Value = -6184.875
$ws.Range("N53").Value = -2029.9
$ws.Range("H64").Value = 7740.4
$ws.Range("I64").Value = 3501
$ws.Range("J64").Value = 10566.667
$ws.Range("K64").Value = 3501
$ws.Range("L64").Value = 10566.667
$ws.Range("M64").Value = -3253
$ws.Range("N64").Value = -11062.667
$ws.Range("H67").Value = 7740.4
$ws.Range("I67").Value = 3501
$ws.Range("J67").Value = 10566.667
$ws.Range("K67").Value = 3501
$ws.Range("L67").Value = 10566.667
$ws.Range("M67").Value = -2643
$ws.Range("N67").Value = -12282.667
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 19416.666
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 19416.666
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -21288.666
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 19416.666
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 97083.33
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -106443.33
$ws.Range("H135").Value = 1170.9
$ws.Range("I135").Value = 1201.1177
$ws.Range("K135").Value = 10810.0593
$ws.Range("M135").Value = -8275.059300000001
$ws.Range("H137").Value = 2187.2727
$ws.Range("I137").Value = 1940.1111
$ws.Range("K137").Value = 5820.3333
$ws.Range("M137").Value = -3270.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 29431
$ws.Range("I31").Value = 21203.4
$ws.Range("K31").Value = 21203.4
$ws.Range("M31").Value = -20909.4
$ws.Range("H45").Value = 5677.0386
$ws.Range("I45").Value = 6337.684
$ws.Range("K45").Value = 6337.684
$ws.Range("M45").Value = -5960.684
$ws.Range("H61").Value = 3869.9783
$ws.Range("I61").Value = 2893.2
$ws.Range("K61").Value = 2893.2
$ws.Range("M61").Value = -2681.2
$ws.Range("H102").Value = 1304.5
$ws.Range("I102").Value = 1221.7778
$ws.Range("J102").Value = 1751.2
$ws.Range("K102").Value = 1221.7778
$ws.Range("L102").Value = 1751.2
$ws.Range("M102").Value = 400.2221999999999
$ws.Range("N102").Value = -4995.2
$ws.Range("H132").Value = 2326
$ws.Range("I132").Value = 2022.6552
$ws.Range("J132").Value = 3125.7273
$ws.Range("K132").Value = 6067.9656
$ws.Range("L132").Value = 9377.1819
$ws.Range("M132").Value = -3537.9656
$ws.Range("N132").Value = -14437.1819
$ws.Range("H136").Value = 3869.9783
$ws.Range("I136").Value = 2893.2
$ws.Range("K136").Value = 8679.599999999999
$ws.Range("M136").Value = -6129.599999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2274.9167
$ws.Range("I86").Value = 3117.3333
$ws.Range("J86").Value = 1432.5
$ws.Range("K86").Value = 3117.3333
$ws.Range("L86").Value = 1432.5
$ws.Range("M86").Value = -1994.3333
$ws.Range("N86").Value = -3678.5
$ws.Range("H88").Value = 32500
$ws.Range("J88").Value = 32500
$ws.Range("L88").Value = 32500
$ws.Range("N88").Value = -33312
$ws.Range("H89").Value = 2274.9167
$ws.Range("I89").Value = 3117.3333
$ws.Range("J89").Value = 1432.5
$ws.Range("K89").Value = 15586.6665
$ws.Range("L89").Value = 7162.5
$ws.Range("M89").Value = -9970.666499999999
$ws.Range("N89").Value = -18394.5
$ws.Range("H91").Value = 32500
$ws.Range("J91").Value = 32500
$ws.Range("L91").Value = 32500
$ws.Range("N91").Value = -35308

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2055.9583
$ws.Range("I31").Value = 1211.625
$ws.Range("J31").Value = 2478.125
$ws.Range("K31").Value = 1211.625
$ws.Range("L31").Value = 2478.125
$ws.Range("M31").Value = -916.625
$ws.Range("N31").Value = -3068.125
$ws.Range("H34").Value = 2055.9583
$ws.Range("I34").Value = 1211.625
$ws.Range("J34").Value = 2478.125
$ws.Range("K34").Value = 1211.625
$ws.Range("L34").Value = 2478.125
$ws.Range("M34").Value = -1009.625
$ws.Range("N34").Value = -2882.125
$ws.Range("H122").Value = 3001.575
$ws.Range("I122").Value = 2716.4666
$ws.Range("K122").Value = 8149.399800000001
$ws.Range("M122").Value = -5699.399800000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1378.7858
$ws.Range("I86").Value = 429.25
$ws.Range("J86").Value = 2644.8333
$ws.Range("K86").Value = 1287.75
$ws.Range("L86").Value = 7934.499899999999
$ws.Range("M86").Value = -101.75
$ws.Range("N86").Value = -10306.4999
$ws.Range("H89").Value = 1378.7858
$ws.Range("I89").Value = 429.25
$ws.Range("J89").Value = 2644.8333
$ws.Range("K89").Value = 3863.25
$ws.Range("L89").Value = 23803.4997
$ws.Range("M89").Value = 2064.75
$ws.Range("N89").Value = -35659.4997
$ws.Range("H139").Value = 2155.1904
$ws.Range("I139").Value = 1967.6
$ws.Range("J139").Value = 2624.1667
$ws.Range("K139").Value = 5902.799999999999
$ws.Range("L139").Value = 7872.500100000001
$ws.Range("M139").Value = -762.7999999999993
$ws.Range("N139").Value = -18152.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4875.607
$ws.Range("I70").Value = 4590.1665
$ws.Range("J70").Value = 5389.4
$ws.Range("K70").Value = 4590.1665
$ws.Range("L70").Value = 5389.4
$ws.Range("M70").Value = -4320.1665
$ws.Range("N70").Value = -5929.4
$ws.Range("H73").Value = 4875.607
$ws.Range("I73").Value = 4590.1665
$ws.Range("J73").Value = 5389.4
$ws.Range("K73").Value = 4590.1665
$ws.Range("L73").Value = 5389.4
$ws.Range("M73").Value = -3654.1665
$ws.Range("N73").Value = -7261.4
$ws.Range("H80").Value = 3799.75
$ws.Range("I80").Value = 2499.5
$ws.Range("J80").Value = 5100
$ws.Range("K80").Value = 2499.5
$ws.Range("L80").Value = 5100
$ws.Range("M80").Value = -1501.5
$ws.Range("N80").Value = -7096
$ws.Range("H83").Value = 3799.75
$ws.Range("I83").Value = 2499.5
$ws.Range("J83").Value = 5100
$ws.Range("K83").Value = 12497.5
$ws.Range("L83").Value = 25500
$ws.Range("M83").Value = -7505.5
$ws.Range("N83").Value = -35484
$ws.Range("H132").Value = 2780.5789
$ws.Range("J132").Value = 3347
$ws.Range("L132").Value = 10041
$ws.Range("N132").Value = -15101

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 53626.5
$ws.Range("I7").Value = 53626.5
$ws.Range("K7").Value = 53626.5
$ws.Range("M7").Value = -53514.5
$ws.Range("H68").Value = 3974.6667
$ws.Range("I68").Value = 2021.7778
$ws.Range("K68").Value = 2021.7778
$ws.Range("M68").Value = -1272.7778
$ws.Range("H71").Value = 3974.6667
$ws.Range("I71").Value = 2021.7778
$ws.Range("K71").Value = 10108.889
$ws.Range("M71").Value = -6364.889000000001
$ws.Range("H122").Value = 4333.702
$ws.Range("I122").Value = 3145.6453
$ws.Range("K122").Value = 9436.9359
$ws.Range("M122").Value = -6986.9359
$ws.Range("H126").Value = 53626.5
$ws.Range("I126").Value = 53626.5
$ws.Range("K126").Value = 160879.5
$ws.Range("M126").Value = -158409.5

Write-Host "Applied 184 cell updates."
